$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: new "semi supervised learning" literature-survey task row.
# D28 already holds the task description (shared string 61); this edit
# fills in start/end "dates", workdays (gross/net) and a comment, matching
# the pattern used by the other rows in the table.

$ws.Range("E28").Value = 44011.447916666664
$ws.Range("F28").Value = "6/30/2020 6/29/20 10:21"
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 1
$ws.Range("I28").Value = "בלוגפוסט שווה אחוששלוקי https://ruder.io/semi-supervised/`nלדעתי, הנושא הזה צריך לבוא אחרי הפרק של רשתות כי יש שיטות שרלוונטיות בעיקר לזה. בכלל, הנושא הזה יותר מתאים לתמונות."
$ws.Range("I28").WrapText = $true

$ws.Rows("28").RowHeight = 43.5

$ws.Range("I29").Select() | Out-Null
